# Nudge the partner/client logo picture rows up very slightly on the two
# "buyer logos" slides (slide 3 and slide 4). Only the vertical position
# (.Top) changes for each picture; horizontal position and size stay put.
#
# NOTE on the literal values below: PowerPoint's Shape.Top/.Left COM
# properties are in points (1 pt = 12700 EMU) and are round-tripped through
# single precision internally, so simply computing target_emu / 12700.0 in
# script can land one EMU short after the runtime's own float32 rounding.
# The constants here are chosen (via the f32 round-trip) so that, after the
# runtime converts them back to EMU, they land exactly on the target EMU
# values from the authoritative OOXML.

$p = $ppt.ActivePresentation

# slide index -> ordered list of new Top values (in points), one per
# picture shape, in the order the pictures appear on the slide.
$targets = @{
    3 = @(107.4574203491211, 167.420654296875, 226.90672302246094, 286.392822265625, 345.8788757324219)
    4 = @(107.4574203491211, 167.420654296875, 226.90672302246094, 286.392822265625)
}

foreach ($slideIndex in $targets.Keys) {
    $slide = $p.Slides.Item($slideIndex)
    $newTops = $targets[$slideIndex]

    $picIndex = 0
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)
        if ($shape.Type -eq 13) {
            # msoPicture
            if ($picIndex -lt $newTops.Count) {
                $shape.Top = $newTops[$picIndex]
            }
            $picIndex = $picIndex + 1
        }
    }
}
